# Auto-generated script applying scheduled-runner profit recalculations
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1185.25
$ws.Range("I28").Value = 840.2727
$ws.Range("J28").Value = 4980
$ws.Range("K28").Value = 840.2727
$ws.Range("L28").Value = 4980
$ws.Range("M28").Value = -355.2727
$ws.Range("N28").Value = -5950
$ws.Range("H32").Value = 927.8
$ws.Range("I32").Value = 897.5
$ws.Range("J32").Value = 935.375
$ws.Range("K32").Value = 897.5
$ws.Range("L32").Value = 935.375
$ws.Range("M32").Value = -571.5
$ws.Range("N32").Value = -1587.375
$ws.Range("H39").Value = 335.57144
$ws.Range("I39").Value = 95.8125
$ws.Range("J39").Value = 1102.8
$ws.Range("K39").Value = 287.4375
$ws.Range("L39").Value = 3308.4
$ws.Range("M39").Value = 8.5625
$ws.Range("N39").Value = -3900.4
$ws.Range("H51").Value = 4964.357
$ws.Range("I51").Value = 1980.2
$ws.Range("J51").Value = 6622.222
$ws.Range("K51").Value = 1980.2
$ws.Range("L51").Value = 6622.222
$ws.Range("M51").Value = -1496.2
$ws.Range("N51").Value = -7590.222
$ws.Range("H132").Value = 2056.9016
$ws.Range("I132").Value = 872.4039
$ws.Range("J132").Value = 8900.666999999999
$ws.Range("K132").Value = 2617.2117
$ws.Range("L132").Value = 26702.001
$ws.Range("M132").Value = -87.21169999999984
$ws.Range("N132").Value = -31762.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 20100
$ws.Range("J17").Value = 20100
$ws.Range("L17").Value = 20100
$ws.Range("N17").Value = -20446
$ws.Range("H21").Value = 6440.273
$ws.Range("I21").Value = 3855.375
$ws.Range("J21").Value = 13333.333
$ws.Range("K21").Value = 3855.375
$ws.Range("L21").Value = 13333.333
$ws.Range("M21").Value = -3481.375
$ws.Range("N21").Value = -14081.333
$ws.Range("H33").Value = 2500
$ws.Range("I33").Value = 2500
$ws.Range("K33").Value = 2500
$ws.Range("M33").Value = -2171
$ws.Range("H45").Value = 2020.1666
$ws.Range("I45").Value = 1733.3334
$ws.Range("J45").Value = 2115.7778
$ws.Range("K45").Value = 1733.3334
$ws.Range("L45").Value = 2115.7778
$ws.Range("M45").Value = -1356.3334
$ws.Range("N45").Value = -2869.7778
$ws.Range("H61").Value = 1575.8723
$ws.Range("I61").Value = 1401.0264
$ws.Range("K61").Value = 1401.0264
$ws.Range("M61").Value = -1189.0264
$ws.Range("H97").Value = 1004.4091
$ws.Range("I97").Value = 918.2632
$ws.Range("K97").Value = 918.2632
$ws.Range("M97").Value = -422.2632
$ws.Range("H102").Value = 1902
$ws.Range("I102").Value = 1627.5
$ws.Range("K102").Value = 1627.5
$ws.Range("M102").Value = -5.5
$ws.Range("H110").Value = 1862
$ws.Range("I110").Value = 1862
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1862
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 183
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 1763.2069
$ws.Range("I122").Value = 1266.0714
$ws.Range("J122").Value = 2227.2
$ws.Range("K122").Value = 3798.2142
$ws.Range("L122").Value = 6681.599999999999
$ws.Range("M122").Value = -1348.2142
$ws.Range("N122").Value = -11581.6
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2037.9672
$ws.Range("I132").Value = 1756.0217
$ws.Range("J132").Value = 2902.6
$ws.Range("K132").Value = 5268.0651
$ws.Range("L132").Value = 8707.799999999999
$ws.Range("M132").Value = -2738.0651
$ws.Range("N132").Value = -13767.8
$ws.Range("H136").Value = 1575.8723
$ws.Range("I136").Value = 1401.0264
$ws.Range("K136").Value = 4203.0792
$ws.Range("M136").Value = -1653.0792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1492.4166
$ws.Range("I94").Value = 1389.9
$ws.Range("J94").Value = 2005
$ws.Range("K94").Value = 1389.9
$ws.Range("L94").Value = 2005
$ws.Range("M94").Value = -938.9000000000001
$ws.Range("N94").Value = -2907
$ws.Range("H97").Value = 10257
$ws.Range("I97").Value = 3676
$ws.Range("J97").Value = 30000
$ws.Range("K97").Value = 3676
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = -2685
$ws.Range("N97").Value = -31982
$ws.Range("H107").Value = 1537.963
$ws.Range("I107").Value = 1448
$ws.Range("J107").Value = 1751.625
$ws.Range("K107").Value = 1448
$ws.Range("L107").Value = 1751.625
$ws.Range("M107").Value = 472
$ws.Range("N107").Value = -5591.625
$ws.Range("H134").Value = 3538.7576
$ws.Range("I134").Value = 2617.8215
$ws.Range("J134").Value = 8696
$ws.Range("K134").Value = 7853.4645
$ws.Range("L134").Value = 26088
$ws.Range("M134").Value = -5318.4645
$ws.Range("N134").Value = -31158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 929.61536
$ws.Range("I58").Value = 798.6957
$ws.Range("J58").Value = 1933.3334
$ws.Range("K58").Value = 798.6957
$ws.Range("L58").Value = 1933.3334
$ws.Range("M58").Value = -595.6957
$ws.Range("N58").Value = -2339.3334
$ws.Range("H122").Value = 1508.8235
$ws.Range("I122").Value = 1366.8572
$ws.Range("J122").Value = 2171.3333
$ws.Range("K122").Value = 4100.571599999999
$ws.Range("L122").Value = 6513.999899999999
$ws.Range("M122").Value = -1650.571599999999
$ws.Range("N122").Value = -11413.9999
$ws.Range("H132").Value = 3889
$ws.Range("I132").Value = 3535.5
$ws.Range("K132").Value = 10606.5
$ws.Range("M132").Value = -8076.5
$ws.Range("H136").Value = 929.61536
$ws.Range("I136").Value = 798.6957
$ws.Range("J136").Value = 1933.3334
$ws.Range("K136").Value = 2396.0871
$ws.Range("L136").Value = 5800.0002
$ws.Range("M136").Value = 153.9129000000003
$ws.Range("N136").Value = -10900.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 764.5714
$ws.Range("J5").Value = 991.2143
$ws.Range("L5").Value = 2973.6429
$ws.Range("N5").Value = -3197.6429
$ws.Range("H131").Value = 809.40405
$ws.Range("I131").Value = 420.26315
$ws.Range("J131").Value = 901.825
$ws.Range("K131").Value = 1260.78945
$ws.Range("L131").Value = 2705.475
$ws.Range("M131").Value = 3779.21055
$ws.Range("N131").Value = -12785.475
$ws.Range("H135").Value = 764.5714
$ws.Range("J135").Value = 991.2143
$ws.Range("L135").Value = 8920.9287
$ws.Range("N135").Value = -13990.9287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352
$ws.Range("H102").Value = 1501.45
$ws.Range("I102").Value = 1286.7333
$ws.Range("J102").Value = 2145.6
$ws.Range("K102").Value = 1286.7333
$ws.Range("L102").Value = 2145.6
$ws.Range("M102").Value = 335.2666999999999
$ws.Range("N102").Value = -5389.6
$ws.Range("H107").Value = 879.95
$ws.Range("I107").Value = 814.1667
$ws.Range("J107").Value = 978.625
$ws.Range("K107").Value = 814.1667
$ws.Range("L107").Value = 978.625
$ws.Range("M107").Value = 1105.8333
$ws.Range("N107").Value = -4818.625
$ws.Range("H122").Value = 27480.05
$ws.Range("I122").Value = 34964.465
$ws.Range("J122").Value = 2532
$ws.Range("K122").Value = 104893.395
$ws.Range("L122").Value = 7596
$ws.Range("M122").Value = -102443.395
$ws.Range("N122").Value = -12496
$ws.Range("H132").Value = 2317.3262
$ws.Range("I132").Value = 2370.4285
$ws.Range("J132").Value = 2148.3635
$ws.Range("K132").Value = 7111.2855
$ws.Range("L132").Value = 6445.0905
$ws.Range("M132").Value = -4581.2855
$ws.Range("N132").Value = -11505.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2333.3333
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 2750
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2750
$ws.Range("M40").Value = -1364
$ws.Range("N40").Value = -3022
$ws.Range("H54").Value = 4984
$ws.Range("J54").Value = 4984
$ws.Range("L54").Value = 4984
$ws.Range("N54").Value = -6272
$ws.Range("H64").Value = 30872.5
$ws.Range("J64").Value = 30872.5
$ws.Range("L64").Value = 30872.5
$ws.Range("N64").Value = -31322.5
$ws.Range("H67").Value = 30872.5
$ws.Range("J67").Value = 30872.5
$ws.Range("L67").Value = 30872.5
$ws.Range("N67").Value = -32432.5
$ws.Range("H99").Value = 28500
$ws.Range("J99").Value = 28500
$ws.Range("L99").Value = 28500
$ws.Range("N99").Value = -34490
$ws.Range("H136").Value = 11495228
$ws.Range("I136").Value = 13889830
$ws.Range("J136").Value = 1141
$ws.Range("K136").Value = 41669490
$ws.Range("L136").Value = 3423
$ws.Range("M136").Value = -41666940
$ws.Range("N136").Value = -8523

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 10810.889
$ws.Range("J69").Value = 10810.889
$ws.Range("L69").Value = 10810.889
$ws.Range("N69").Value = -12308.889
$ws.Range("H72").Value = 10810.889
$ws.Range("J72").Value = 10810.889
$ws.Range("L72").Value = 32432.667
$ws.Range("N72").Value = -39920.667
$ws.Range("H113").Value = 488.2143
$ws.Range("I113").Value = 424.83334
$ws.Range("J113").Value = 535.75
$ws.Range("K113").Value = 1274.50002
$ws.Range("L113").Value = 1607.25
$ws.Range("M113").Value = 895.4999800000001
$ws.Range("N113").Value = -5947.25
$ws.Range("H122").Value = 8973171
$ws.Range("I122").Value = 10871361
$ws.Range("K122").Value = 32614083
$ws.Range("M122").Value = -32611633
$ws.Range("H132").Value = 3127306
$ws.Range("I132").Value = 3848353.8
$ws.Range("J132").Value = 2765.5
$ws.Range("K132").Value = 11545061.4
$ws.Range("L132").Value = 8296.5
$ws.Range("M132").Value = -11542531.4
$ws.Range("N132").Value = -13356.5
$ws.Range("H136").Value = 27780274
$ws.Range("I136").Value = 41669496
$ws.Range("K136").Value = 125008488
$ws.Range("M136").Value = -125005938
